# Regenerate merged AHB files
#
# 1. Rename the "_old" / "_new" suffixed header columns to "_FV2310" / "_FV2404"
# 2. Turn the data range A1:U64 into an Excel Table ("Table1") with an AutoFilter
# 3. Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) ----------------------------------------
$fv2310 = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

$fv2404 = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $fv2310.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2310[$i]
}
# column 11 ("diff") keeps its name
for ($i = 0; $i -lt $fv2404.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2404[$i]
}

# --- 2. Convert the data range into an Excel Table ------------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U64"), 0, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ---------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
